$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 232, shifting rows 232:289 down to 233:290
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with new data
$ws.Cells.Item(232, 1).Value = 5
$ws.Cells.Item(232, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(232, 3).Value = "Maule"
$ws.Cells.Item(232, 4).Value = 44951
$ws.Cells.Item(232, 4).NumberFormat = $ws.Cells.Item(233, 4).NumberFormat
$ws.Cells.Item(232, 5).Value = 7
$ws.Cells.Item(232, 6).Value = 100112021
$ws.Cells.Item(232, 7).Value = "Ají"
$ws.Cells.Item(232, 8).Value = "Americana (o)"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 250
$ws.Cells.Item(232, 11).Value = 9000
$ws.Cells.Item(232, 12).Value = 10000
$ws.Cells.Item(232, 13).Value = 9400
$ws.Cells.Item(232, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(232, 15).Value = "Región del Maule"
$ws.Cells.Item(232, 16).Value = 627
$ws.Cells.Item(232, 17).Value = 15
$ws.Cells.Item(232, 18).Value = "Hortaliza"
